$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 85.30768999999999
$ws.Range("I39").Value = 85.30768999999999
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 255.92307
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = 40.07693
$ws.Range("N39").ClearContents()
$ws.Range("H43").Value = 3983.8
$ws.Range("I43").Value = 2290
$ws.Range("J43").Value = 5677.6
$ws.Range("K43").Value = 2290
$ws.Range("L43").Value = 5677.6
$ws.Range("M43").Value = -2221
$ws.Range("N43").Value = -5815.6
$ws.Range("H96").Value = 3974.2856
$ws.Range("I96").Value = 3200
$ws.Range("J96").Value = 5006.6665
$ws.Range("K96").Value = 9600
$ws.Range("L96").Value = 15019.9995
$ws.Range("M96").Value = -8227
$ws.Range("N96").Value = -17765.9995
$ws.Range("H98").Value = 838.5
$ws.Range("J98").Value = 600
$ws.Range("L98").Value = 600
$ws.Range("N98").Value = -3596
$ws.Range("H107").Value = 396.14285
$ws.Range("I107").Value = 396.14285
$ws.Range("K107").Value = 396.14285
$ws.Range("M107").Value = 1523.85715
$ws.Range("H111").Value = 3956.3333
$ws.Range("I111").Value = 3956.3333
$ws.Range("K111").Value = 11868.9999
$ws.Range("M111").Value = -8801.999899999999
$ws.Range("H113").Value = 3484.8462
$ws.Range("I113").Value = 3029.1428
$ws.Range("K113").Value = 3029.1428
$ws.Range("M113").Value = 224.8571999999999
$ws.Range("H122").Value = 838.5
$ws.Range("J122").Value = 600
$ws.Range("L122").Value = 1800
$ws.Range("N122").Value = -6700
$ws.Range("H125").Value = 2359.75
$ws.Range("I125").Value = 969.5
$ws.Range("J125").Value = 3750
$ws.Range("K125").Value = 8725.5
$ws.Range("L125").Value = 33750
$ws.Range("M125").Value = -6265.5
$ws.Range("N125").Value = -38670

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1829.2858
$ws.Range("I61").Value = 1829.2858
$ws.Range("K61").Value = 1829.2858
$ws.Range("M61").Value = -1617.2858
$ws.Range("H97").Value = 2061.2222
$ws.Range("I97").Value = 471.5
$ws.Range("J97").Value = 5240.6665
$ws.Range("K97").Value = 471.5
$ws.Range("L97").Value = 5240.6665
$ws.Range("M97").Value = 24.5
$ws.Range("N97").Value = -6232.6665
$ws.Range("H136").Value = 1829.2858
$ws.Range("I136").Value = 1829.2858
$ws.Range("K136").Value = 5487.857400000001
$ws.Range("M136").Value = -2937.857400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3870.4119
$ws.Range("I86").Value = 3796.75
$ws.Range("J86").Value = 4047.2
$ws.Range("K86").Value = 3796.75
$ws.Range("L86").Value = 4047.2
$ws.Range("M86").Value = -2673.75
$ws.Range("N86").Value = -6293.2
$ws.Range("H89").Value = 3870.4119
$ws.Range("I89").Value = 3796.75
$ws.Range("J89").Value = 4047.2
$ws.Range("K89").Value = 18983.75
$ws.Range("L89").Value = 20236
$ws.Range("M89").Value = -13367.75
$ws.Range("N89").Value = -31468
$ws.Range("H94").Value = 3900
$ws.Range("I94").Value = 3866.6667
$ws.Range("K94").Value = 3866.6667
$ws.Range("M94").Value = -3415.6667
$ws.Range("H99").Value = 2675.0386
$ws.Range("I99").Value = 2437.2173
$ws.Range("K99").Value = 2437.2173
$ws.Range("M99").Value = -939.2172999999998
$ws.Range("H105").Value = 4464.45
$ws.Range("I105").Value = 3788.7856
$ws.Range("K105").Value = 3788.7856
$ws.Range("M105").Value = -2041.7856
$ws.Range("H122").Value = 299988.66
$ws.Range("J122").Value = 299990
$ws.Range("L122").Value = 299990
$ws.Range("N122").Value = -309790

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 15345.053
$ws.Range("I99").Value = 13900.667
$ws.Range("K99").Value = 13900.667
$ws.Range("M99").Value = -12402.667
$ws.Range("H121").Value = 29997.5
$ws.Range("J121").Value = 29997.5
$ws.Range("L121").Value = 29997.5
$ws.Range("N121").Value = -32617.5
$ws.Range("H126").Value = 15345.053
$ws.Range("I126").Value = 13900.667
$ws.Range("K126").Value = 41702.001
$ws.Range("M126").Value = -39232.001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 58894.41
$ws.Range("J2").Value = 142.66667
$ws.Range("L2").Value = 856.0000200000001
$ws.Range("N2").Value = -1082.00002
$ws.Range("H39").Value = 808.9
$ws.Range("I39").Value = 343.33334
$ws.Range("J39").Value = 4999
$ws.Range("K39").Value = 1030.00002
$ws.Range("L39").Value = 14997
$ws.Range("M39").Value = -736.0000199999999
$ws.Range("N39").Value = -15585
$ws.Range("H48").Value = 512.5
$ws.Range("I48").Value = 450
$ws.Range("J48").Value = 533.3333
$ws.Range("K48").Value = 1350
$ws.Range("L48").Value = 1599.9999
$ws.Range("M48").Value = -1100
$ws.Range("N48").Value = -2099.9999
$ws.Range("H122").Value = 727.3333
$ws.Range("I122").Value = 608.8
$ws.Range("J122").Value = 875.5
$ws.Range("K122").Value = 5479.2
$ws.Range("L122").Value = 7879.5
$ws.Range("M122").Value = -3029.2
$ws.Range("N122").Value = -12779.5
$ws.Range("H136").Value = 13503.143
$ws.Range("I136").Value = 6904.6
$ws.Range("J136").Value = 29999.5
$ws.Range("K136").Value = 20713.8
$ws.Range("L136").Value = 89998.5
$ws.Range("M136").Value = -15613.8
$ws.Range("N136").Value = -100198.5
$ws.Range("H141").Value = 3549.75
$ws.Range("I141").Value = 2233.1667
$ws.Range("K141").Value = 6699.500100000001
$ws.Range("M141").Value = -1519.500100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5450
$ws.Range("J80").Value = 10000
$ws.Range("L80").Value = 10000
$ws.Range("N80").Value = -11996
$ws.Range("H83").Value = 5450
$ws.Range("J83").Value = 10000
$ws.Range("L83").Value = 50000
$ws.Range("N83").Value = -59984
$ws.Range("H97").Value = 1371.5714
$ws.Range("J97").Value = 1439.1666
$ws.Range("L97").Value = 1439.1666
$ws.Range("N97").Value = -2431.1666
$ws.Range("H113").Value = 4985.7144
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H122").Value = 204151.8
$ws.Range("I122").Value = 4327.6665
$ws.Range("J122").Value = 503888
$ws.Range("K122").Value = 12982.9995
$ws.Range("L122").Value = 1511664
$ws.Range("M122").Value = -10532.9995
$ws.Range("N122").Value = -1516564
$ws.Range("H126").Value = 4147.6665
$ws.Range("I126").Value = 2510
$ws.Range("J126").Value = 4966.5
$ws.Range("K126").Value = 7530
$ws.Range("L126").Value = 14899.5
$ws.Range("M126").Value = -5060
$ws.Range("N126").Value = -19839.5
$ws.Range("H132").Value = 2784.125
$ws.Range("I132").Value = 1973.7693
$ws.Range("J132").Value = 6295.6665
$ws.Range("K132").Value = 5921.3079
$ws.Range("L132").Value = 18886.9995
$ws.Range("M132").Value = -3391.3079
$ws.Range("N132").Value = -23946.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2795.4
$ws.Range("I7").Value = 2665.6667
$ws.Range("K7").Value = 2665.6667
$ws.Range("M7").Value = -2553.6667
$ws.Range("H22").Value = 4249.75
$ws.Range("I22").Value = 3999.6667
$ws.Range("K22").Value = 3999.6667
$ws.Range("M22").Value = -3704.6667
$ws.Range("H27").Value = 4249.75
$ws.Range("I27").Value = 3999.6667
$ws.Range("K27").Value = 3999.6667
$ws.Range("M27").Value = -3892.6667
$ws.Range("H56").Value = 12030.8
$ws.Range("I56").Value = 10013.667
$ws.Range("K56").Value = 10013.667
$ws.Range("M56").Value = -9322.666999999999
$ws.Range("H126").Value = 2795.4
$ws.Range("I126").Value = 2665.6667
$ws.Range("K126").Value = 7997.000100000001
$ws.Range("M126").Value = -5527.000100000001
$ws.Range("H132").Value = 3439.7856
$ws.Range("I132").Value = 2906.3428
$ws.Range("J132").Value = 6107
$ws.Range("K132").Value = 8719.028399999999
$ws.Range("L132").Value = 18321
$ws.Range("M132").Value = -6189.028399999999
$ws.Range("N132").Value = -23381
$ws.Range("H136").Value = 1930.8182
$ws.Range("I136").Value = 1950
$ws.Range("K136").Value = 5850
$ws.Range("M136").Value = -3300

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 510000
$ws.Range("J49").Value = 20000
$ws.Range("L49").Value = 20000
$ws.Range("N49").Value = -20460
$ws.Range("H126").Value = 7900
$ws.Range("I126").Value = 800
$ws.Range("J126").Value = 15000
$ws.Range("K126").Value = 2400
$ws.Range("L126").Value = 45000
$ws.Range("M126").Value = 70
$ws.Range("N126").Value = -49940
$ws.Range("H136").Value = 4652.2856
$ws.Range("I136").Value = 1230.7142
$ws.Range("K136").Value = 3692.1426
$ws.Range("M136").Value = -1142.1426
